$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 6.767081260681152
$ws.Range("B1").Value = 5.657371044158936
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 3.16493821144104
$ws.Range("E1").Value = 1.867285370826721
